$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Insert three new slides right after the title/splash slide (position 1),
# i.e. at positions 2, 3 and 4, pushing the existing "First Round of
# Investment" / "Azure Deployment" / "Done" slides down to positions 5-7.
# Layout 2 == "Title and Content" (same layout used by the other slides).
# ---------------------------------------------------------------------------

# --- New slide 2: "13:00" -------------------------------------------------
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "13:00"

# --- New slide 3: "22:00" -------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "22:00"

# --- New slide 4: "Thank You!" --------------------------------------------
$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Thank You!"

$body = $s4.Shapes.Item(2).TextFrame.TextRange
$body.Text = "https://github.com/Isantipov/HotScroll`rTeam:  `rivan.antsipau@gmail.com`rd.parf@live.com`re.leychenok@gmail.com`r`r19graff91@gmail.com   `r`r"

# Paragraph indent levels (1 = top level, 2 = first sub-level -> pPr lvl="1")
for ($i = 3; $i -le 8; $i++) {
    $body.Paragraphs($i, 1).IndentLevel = 2
}

# Hyperlink on the project URL - split the same way AutoFormat would (the
# "https://" prefix and the rest of the address as separate runs sharing the
# same hyperlink relationship).
$urlPart1 = $body.Characters(1, 8)
$urlPart1.ActionSettings(1).Hyperlink.Address = "https://github.com/Isantipov/HotScroll"
$urlPart2 = $body.Characters(9, 30)
$urlPart2.ActionSettings(1).Hyperlink.Address = "https://github.com/Isantipov/HotScroll"

# Team member e-mails (first three share the same mailto link, matching the
# source deck; the last one points to a different address).
$email1 = $body.Characters(48, 23)
$email1.ActionSettings(1).Hyperlink.Address = "mailto:ivan.antsipau@gmail.com"

$email2 = $body.Characters(72, 15)
$email2.ActionSettings(1).Hyperlink.Address = "mailto:ivan.antsipau@gmail.com"

$email3 = $body.Characters(88, 21)
$email3.ActionSettings(1).Hyperlink.Address = "mailto:ivan.antsipau@gmail.com"

$email4 = $body.Characters(111, 19)
$email4.ActionSettings(1).Hyperlink.Address = "mailto:19graff91@gmail.com"

Write-Output $p.Slides.Count
